# Atualização de bases das ligas, do dia: 21-02-2024 às 23:25
#
# The underlying data rows got re-sorted (their natural match order changed),
# which manifests in the OOXML as pairs (and one 3-cycle) of data rows whose
# content - everything except the running index (col A) and the constant
# Div / Div Original Name / Date columns (C, D, E) - gets rotated among the
# rows that share the same Date (E) value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that participate in the row rotation (1-based column indices):
# B=2 (id), F=6..AC=29 (everything from HomeTeam onward)
$swapCols = @(2,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,25,26,27,28,29)

# Groups of row numbers whose content rotates: new(group[i]) = old(group[i+1 mod n])
$groups = @(
    ,@(68,69)
    ,@(79,80)
    ,@(98,99)
    ,@(105,106)
    ,@(141,142)
    ,@(193,194)
    ,@(220,221)
    ,@(274,275)
    ,@(296,297)
    ,@(298,299)
    ,@(314,315)
    ,@(386,387)
    ,@(394,395,396)
    ,@(397,398)
)

foreach ($group in $groups) {
    $n = $group.Count

    # Snapshot the current ("old") values of every swap column for every row
    # in this group before writing anything, so the rotation doesn't clobber
    # data it still needs to read.
    $snapshot = @{}
    for ($i = 0; $i -lt $n; $i++) {
        $row = $group[$i]
        $rowVals = @{}
        foreach ($col in $swapCols) {
            $rowVals[$col] = $ws.Cells.Item($row, $col).Value2
        }
        $snapshot[$row] = $rowVals
    }

    # Write back: row at position i receives the old content of the row at
    # position (i+1) mod n, i.e. a left-rotation of the group.
    for ($i = 0; $i -lt $n; $i++) {
        $destRow = $group[$i]
        $srcRow = $group[($i + 1) % $n]
        $srcVals = $snapshot[$srcRow]
        foreach ($col in $swapCols) {
            $ws.Cells.Item($destRow, $col).Value = $srcVals[$col]
        }
    }
}
